$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.237.57'
$ws.Range('E2').Value = '  +5.29%  '

$ws.Range('D3').Value = '3.513.07'
$ws.Range('E3').Value = '  +2.83%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '595.03'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +4.29%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '169.54'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +7.24%  '

$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '1.00'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  -0.06%  '

$ws.Range('D8').Value = '3.511.67'
$ws.Range('E8').Value = '  +2.73%  '

$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.576'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  +1.81%  '

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '7.29'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +0.75%  '

$ws.Range('E11').Value = '  +5.87%  '

$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.442'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +4.96%  '

$ws.Range('D13').Value = '4.115.52'
$ws.Range('E13').Value = '  +2.84%  '

$ws.Range('E14').Value = '  +0.21%  '

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '28.36'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +4.84%  '

$ws.Range('E16').Value = '  +4.63%  '

$ws.Range('D17').Value = '67.150.05'
$ws.Range('E17').Value = '  +5.04%  '

$ws.Range('D18').Value = '3.508.90'
$ws.Range('E18').Value = '  +2.77%  '

$ws.Range('E19').Value = '  +4.45%  '

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '14.07'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +3.42%  '

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '396.06'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +3.28%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '7.99'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +2.63%  '

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '73.48'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +3.00%  '

$ws.Range('E24').Value = '  +11.78%  '

$ws.Range('E25').Value = '  -0.05%  '

$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '0.532'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +3.41%  '

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '10.20'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +5.59%  '

$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '0.183'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +2.49%  '

$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '0.999'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  +0.01%  '

$ws.Range('E30').Value = '  +5.61%  '

$ws.Range('E31').Value = '  +6.74%  '

$ws.Range('E32').Value = '  +4.40%  '

$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '23.65'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +3.44%  '

$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '7.48'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +7.67%  '

$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.00'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +0.08%  '

$ws.Range('E36').Value = '  +6.71%  '

$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '162.12'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  +0.84%  '

$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.901'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +6.48%  '

$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '1.95'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +6.91%  '

$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.0755'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +4.45%  '

$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '4.68'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +7.49%  '

$ws.Range('B42').Value = 'RenderToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '6.73'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +5.42%  '

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '26.56'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +2.36%  '

$ws.Range('B44').Value = 'Maker'
$ws.Range('C44').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D44').Value = '2.852.54'
$ws.Range('E44').Value = '  +2.10%  '

$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '43.51'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +1.01%  '

$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '26.55'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +0.87%  '

$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '2.57'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +7.32%  '

$ws.Range('E48').Value = '  +4.28%  '

$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '352.15'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  +6.43%  '

$ws.Range('E50').Value = '  +5.01%  '

$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '33.64'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +12.59%  '
